# "added interactive graphs and solver"
# The edge list in Sheet1 (A:B = From node / To node) is rebuilt with new
# values, and the trailing rows (11-18) that used to hold the last batch of
# edges are cleared out (values removed, thin-box border removed) to make
# room for the interactive graph/solver area that now lives below/around
# the table. The view is also reset to 100% zoom with the selection parked
# on F9 (inside that new working area).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing edge rows (2-10) with their new From/To values -------
$ws.Range("B2").Value  = 3

$ws.Range("B4").Value  = 6

$ws.Range("A5").Value  = 4
$ws.Range("B5").Value  = 6

$ws.Range("B6").Value  = 7

$ws.Range("A7").Value  = 6

$ws.Range("A8").Value  = 4
$ws.Range("B8").Value  = 9

$ws.Range("A9").Value  = 6
$ws.Range("B9").Value  = 8

$ws.Range("A10").Value = 7
$ws.Range("B10").Value = 8

# --- Rows 11-18 used to hold the rest of the edge list; clear them out and
#     drop their border so they fall back to a plain (unboxed) style ------
$tail = $ws.Range("A11:B18")
$tail.Borders.LineStyle = -4142
$tail.Style = "Normal"
$tail.ClearContents()

# --- Reset the view: 100% zoom, selection on F9 ----------------------------
$excel.ActiveWindow.Zoom = 100
$ws.Range("F9").Select() | Out-Null
